$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2528.8923
$ws.Range("J138").Value = 2129.6833
$ws.Range("L138").Value = 6389.0499
$ws.Range("N138").Value = -16669.0499

$ws.Range("H141").Value = 2650.36
$ws.Range("I141").Value = 885.825
$ws.Range("J141").Value = 9708.5
$ws.Range("K141").Value = 2657.475
$ws.Range("L141").Value = 29125.5
$ws.Range("M141").Value = 2522.525
$ws.Range("N141").Value = -39485.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11472.417
$ws.Range("I28").Value = 3559.0908
$ws.Range("J28").Value = 98519
$ws.Range("K28").Value = 3559.0908
$ws.Range("L28").Value = 98519
$ws.Range("M28").Value = -3367.0908
$ws.Range("N28").Value = -98903

$ws.Range("H68").Value = 95000
$ws.Range("J68").Value = 95000
$ws.Range("L68").Value = 95000
$ws.Range("N68").Value = -96622

$ws.Range("H71").Value = 95000
$ws.Range("J71").Value = 95000
$ws.Range("L71").Value = 285000
$ws.Range("N71").Value = -293112

$ws.Range("H99").Value = 11472.417
$ws.Range("I99").Value = 3559.0908
$ws.Range("J99").Value = 98519
$ws.Range("K99").Value = 3559.0908
$ws.Range("L99").Value = 98519
$ws.Range("M99").Value = -564.0907999999999
$ws.Range("N99").Value = -104509

$ws.Range("H108").Value = 98684
$ws.Range("J108").Value = 98684
$ws.Range("L108").Value = 98684
$ws.Range("N108").Value = -106364

$ws.Range("H129").Value = 38709.89
$ws.Range("J129").Value = 38709.89
$ws.Range("L129").Value = 38709.89
$ws.Range("N129").Value = -48709.89

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H50").Value = 53390
$ws.Range("J50").Value = 53390
$ws.Range("L50").Value = 53390
$ws.Range("N50").Value = -54538

$ws.Range("H96").Value = 5206.75
$ws.Range("I96").Value = 5206.75
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 5206.75
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2460.75
$ws.Range("N96").ClearContents()

$ws.Range("H97").Value = 24966.334
$ws.Range("I97").Value = 4214
$ws.Range("J97").Value = 66471
$ws.Range("K97").Value = 4214
$ws.Range("L97").Value = 66471
$ws.Range("M97").Value = -3223
$ws.Range("N97").Value = -68453

$ws.Range("H109").Value = 49434
$ws.Range("J109").Value = 49434
$ws.Range("L109").Value = 49434
$ws.Range("N109").Value = -52208

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H115").Value = 64228
$ws.Range("J115").Value = 64228
$ws.Range("L115").Value = 64228
$ws.Range("N115").Value = -67362

$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178

$ws.Range("H118").Value = 45387.43
$ws.Range("J118").Value = 45387.43
$ws.Range("L118").Value = 45387.43
$ws.Range("N118").Value = -48701.43

$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 49995
$ws.Range("J53").Value = 49995
$ws.Range("L53").Value = 49995
$ws.Range("N53").Value = -51209

$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

$ws.Range("H112").Value = 38888.25
$ws.Range("J112").Value = 38888.25
$ws.Range("L112").Value = 38888.25
$ws.Range("N112").Value = -41842.25

$ws.Range("H114").Value = 43280.668
$ws.Range("J114").Value = 43280.668
$ws.Range("L114").Value = 43280.668
$ws.Range("N114").Value = -51958.668

$ws.Range("H119").Value = 37748.75
$ws.Range("J119").Value = 37748.75
$ws.Range("L119").Value = 37748.75
$ws.Range("N119").Value = -47424.75

$ws.Range("H127").Value = 56260
$ws.Range("J127").Value = 56260
$ws.Range("L127").Value = 56260
$ws.Range("N127").Value = -66180

$ws.Range("H131").Value = 36442
$ws.Range("J131").Value = 36442
$ws.Range("L131").Value = 36442
$ws.Range("N131").Value = -46522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5153.6924
$ws.Range("J62").Value = 5153.6924
$ws.Range("L62").Value = 15461.0772
$ws.Range("N62").Value = -16833.0772

$ws.Range("H65").Value = 5153.6924
$ws.Range("J65").Value = 5153.6924
$ws.Range("L65").Value = 46383.2316
$ws.Range("N65").Value = -53247.2316

$ws.Range("H87").Value = 2866.6667
$ws.Range("I87").Value = 800
$ws.Range("K87").Value = 2400
$ws.Range("M87").Value = -1152

$ws.Range("H90").Value = 2866.6667
$ws.Range("I90").Value = 800
$ws.Range("K90").Value = 7200
$ws.Range("M90").Value = -960

$ws.Range("H131").Value = 3133.7273
$ws.Range("J131").Value = 3520.625
$ws.Range("L131").Value = 10561.875
$ws.Range("N131").Value = -20641.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28630.5
$ws.Range("J39").Value = 28630.5
$ws.Range("L39").Value = 28630.5
$ws.Range("N39").Value = -29694.5

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

$ws.Range("H92").Value = 25537.75
$ws.Range("J92").Value = 25537.75
$ws.Range("L92").Value = 25537.75
$ws.Range("N92").Value = -29281.75

$ws.Range("H99").Value = 4830.7617
$ws.Range("I99").Value = 3029.7334
$ws.Range("K99").Value = 3029.7334
$ws.Range("M99").Value = -783.7334000000001

$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws.Range("H130").Value = 60080
$ws.Range("J130").Value = 60080
$ws.Range("L130").Value = 60080
$ws.Range("N130").Value = -70120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6195.4443
$ws.Range("I7").Value = 7001.3335
$ws.Range("J7").Value = 5792.5
$ws.Range("K7").Value = 7001.3335
$ws.Range("L7").Value = 5792.5
$ws.Range("M7").Value = -6889.3335
$ws.Range("N7").Value = -6016.5

$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31498

$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -97488

$ws.Range("H70").Value = 85000
$ws.Range("J70").Value = 85000
$ws.Range("L70").Value = 85000
$ws.Range("N70").Value = -85540

$ws.Range("H73").Value = 85000
$ws.Range("J73").Value = 85000
$ws.Range("L73").Value = 85000
$ws.Range("N73").Value = -86872

$ws.Range("H126").Value = 6195.4443
$ws.Range("I126").Value = 7001.3335
$ws.Range("J126").Value = 5792.5
$ws.Range("K126").Value = 21004.0005
$ws.Range("L126").Value = 17377.5
$ws.Range("M126").Value = -18534.0005
$ws.Range("N126").Value = -22317.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1080
$ws.Range("I126").Value = 1080
$ws.Range("K126").Value = 3240
$ws.Range("M126").Value = -770

$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
